$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.445.56'
$ws.Range("E2").Value = '  +7.21%  '

$ws.Range("D3").Value = '2.381.08'
$ws.Range("E3").Value = '  +4.36%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").Value = "'113.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.96%  '

$ws.Range("D6").Value = "'317.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.01%  '

$ws.Range("D7").Value = "'0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.81%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = "'0.627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.51%  '

$ws.Range("E10").Value = '  +10.87%  '

$ws.Range("D11").Value = "'0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.55%  '

$ws.Range("D12").Value = "'8.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.22%  '

$ws.Range("E13").Value = '  +1.59%  '

$ws.Range("D14").Value = "'1.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.39%  '

$ws.Range("D15").Value = "'15.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.15%  '

$ws.Range("D16").Value = '2.742.42'
$ws.Range("E16").Value = '  +4.50%  '

$ws.Range("D17").Value = '2.367.15'
$ws.Range("E17").Value = '  +3.63%  '

$ws.Range("D18").Value = '45.301.75'
$ws.Range("E18").Value = '  +6.43%  '

$ws.Range("D19").Value = "'7.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.01%  '

$ws.Range("E20").Value = '  +3.56%  '

$ws.Range("D21").Value = "'13.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.59%  '

$ws.Range("D22").Value = "'74.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.89%  '

$ws.Range("D23").Value = "'3.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.55%  '

$ws.Range("D24").Value = "'269.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.26%  '

$ws.Range("D25").Value = "'2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.09%  '

$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("D27").Value = "'7.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.08%  '

$ws.Range("D28").Value = "'11.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.51%  '

$ws.Range("E29").Value = '  +1.95%  '

$ws.Range("D30").Value = "'39.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.74%  '

$ws.Range("D31").Value = "'22.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.08%  '

$ws.Range("D32").Value = "'0.0961"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +13.36%  '

$ws.Range("D33").Value = "'171.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.22%  '

$ws.Range("E34").Value = '  +15.54%  '

$ws.Range("E35").Value = '  +2.72%  '

$ws.Range("E36").Value = '  +8.21%  '

$ws.Range("D37").Value = "'4.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.28%  '

$ws.Range("D38").Value = "'3.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.80%  '

$ws.Range("D39").Value = "'4.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.25%  '

$ws.Range("E40").Value = '  +5.88%  '

$ws.Range("E41").Value = '  +10.73%  '

$ws.Range("D42").Value = "'105.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.90%  '

$ws.Range("E43").Value = '  +6.79%  '

$ws.Range("D44").Value = "'71.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.62%  '

$ws.Range("D45").Value = "'13.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.55%  '

$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").Value = "'5.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +13.09%  '

$ws.Range("D48").Value = "'116.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.20%  '

$ws.Range("D49").Value = "'1.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +19.77%  '

$ws.Range("E50").Value = '  +8.36%  '

$ws.Range("D51").Value = "'79.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.11%  '
